$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "Átlagos Lakásárak"

# Data: Year, B (avg m2 price), C (avg lakás méret m2), D (avg lakás ár Ft)
$data = @(
    @(2010, 199.3333333333333, 70, 13953333.33333333),
    @(2011, 192.6666666666667, 70, 13486666.66666667),
    @(2012, 184.6666666666667, 70, 12926666.66666667),
    @(2013, 174.6666666666667, 70, 12226666.66666667),
    @(2014, 178.3333333333333, 70, 12483333.33333333),
    @(2015, 203.3333333333333, 70, 14233333.33333333),
    @(2016, 237.3333333333333, 70, 16613333.33333334),
    @(2017, 268.6666666666667, 70, 18806666.66666667),
    @(2018, 318.6666666666667, 70, 22306666.66666667),
    @(2019, 389.3333333333333, 70, 27253333.33333333),
    @(2020, 413.6666666666667, 70, 28956666.66666667),
    @(2021, 481, 70, 33670000),
    @(2022, 583, 70, 40810000),
    @(2023, 595, 70, 41650000)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
